# Weighting & Scaling update & heatmap
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Scaling sheet: add "Optimal" / "Threshold" columns (E, F), copying D1's
# number/font formatting so no new cell styles are minted ---
$ws2.Range("D1").Copy()
$ws2.Range("E1:F1").PasteSpecial(-4122)
$ws2.Range("E1").Value = "Optimal"
$ws2.Range("F1").Value = "Threshold"

# --- Update Min/Max values for e_modulus, tensile_strain_at_break,
# tensile_yield_strength rows ---
$ws2.Range("B2").Value = 2404
$ws2.Range("C2").Value = 2884.93

$ws2.Range("B3").Value = 1.6
$ws2.Range("C3").Value = 1.9

$ws2.Range("B4").Value = 33
$ws2.Range("C4").Value = 45.12

# --- Restore the active selections on each sheet ---
[void]$ws1.Range("C17").Select()
[void]$ws2.Range("D10").Select()
[void]$ws2.Activate()

Write-Output "done"
